$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: column B now represents "dayStage" instead of "PrescriptionScheduleEntry"
$ws.Range("B1").Value = "dayStage"

# Map old numeric PrescriptionScheduleEntry codes to the new DAYSTAGE enum text.
$map = @{ 2 = "MORNING"; 3 = "MIDDAY"; 4 = "AFTERNOON" }

# Capture the existing codes first (values are being overwritten in-place).
$codes = @{}
for ($r = 2; $r -le 25; $r++) {
    $codes[$r] = [int]$ws.Cells.Item($r, 2).Value2
}

# Write MORNING/AFTERNOON/MIDDAY rows in that pass order so new shared
# strings are registered as MORNING, AFTERNOON, MIDDAY.
foreach ($code in 2, 4, 3) {
    for ($r = 2; $r -le 25; $r++) {
        if ($codes[$r] -eq $code) {
            $ws.Cells.Item($r, 2).Value = $map[$code]
        }
    }
}

# Widen column B to fit the new text values
$ws.Columns.Item(2).ColumnWidth = 21.6

# Update the active selection to reflect the edited range
$ws.Range("B2:B25").Select()
